$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title at the top of the document.
$metaOld = "Meta description: Read the review of ARRR! 10K Ways slot game - play for free and enjoy a low volatility slot game with excellent graphics and cascade mechanism."
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $metaOld) {
        $p.Range.Delete()
        break
    }
}

# 2. Insert a new bold-text paragraph "Play ARRR! 10K Ways Free - Low
#    Volatility Slot Game" right before the final (image-prompt) paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play ARRR! 10K Ways Free - Low Volatility Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($titleXml)

# 3. Replace the image-generation-prompt text (now the last paragraph) with
#    the meta-description copy, keeping its existing (italic) formatting.
$oldPrompt = "Create a feature image for ARRR! 10K Ways that captures the adventurous spirit of a pirate-themed game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing on a sandy beach with palm trees in the background, holding a treasure chest filled with gold coins and jewels. The warrior's clothing and accessories should suggest that they are a pirate on their quest for treasure. The image should also include the game title, ARRR! 10K Ways, in bold and eye-catching letters. The overall design should be bright, colorful, and playful to attract players' attention and generate excitement about the game."
$newPrompt = "Read the review of ARRR! 10K Ways slot game - play for free and enjoy a low volatility slot game with excellent graphics and cascade mechanism."
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newPrompt, 2)
